# Oregon CGT validator mock data fixes
# - Regenerate provider TINs (9-digit) on the PROV_ID sheet and everywhere
#   they are cross-referenced (RX_MED_PROV), keeping provider names aligned.
# - Clear stray "SUB IPA 1" sub-IPA labels that shouldn't apply to
#   Bend Medical Center / Coastal Health Alliance (PROV_ID, TME_PROV).
# - Refresh the sample TME_ALL aggregate numbers for LOB rows 1 and 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 9. PROV_ID — new TINs + clear two incorrect Sub-IPA labels
# ---------------------------------------------------------------------------
$provId = $wb.Worksheets.Item("9. PROV_ID")

$provId.Range("C10").Value = "786579303"
$provId.Range("C11").Value = "896233790"
$provId.Range("C12").Value = "339670711"
$provId.Range("C13").Value = "210053353"

$provId.Range("B14").Value = ""
$provId.Range("C14").Value = "685582861"

$provId.Range("C15").Value = "553035110"
$provId.Range("C16").Value = "200604502"
$provId.Range("C17").Value = "642621108"

$provId.Range("B18").Value = ""
$provId.Range("C18").Value = "702632297"

$provId.Range("C19").Value = "797808098"

# ---------------------------------------------------------------------------
# 2. TME_ALL — updated LOB 1 and LOB 2 aggregate figures
# ---------------------------------------------------------------------------
$tmeAll = $wb.Worksheets.Item("2. TME_ALL")

$tmeAll.Range("C10").Value = 32493
$tmeAll.Range("D10").Value = 0.888
$tmeAll.Range("E10").Value = 10328001.97
$tmeAll.Range("F10").Value = 8819073.57
$tmeAll.Range("G10").Value = 1635208.21
$tmeAll.Range("H10").Value = 5867648.49
$tmeAll.Range("I10").Value = 1330399.3
$tmeAll.Range("J10").Value = 1306462.8

$tmeAll.Range("C11").Value = 15189
$tmeAll.Range("D11").Value = 0.886
$tmeAll.Range("E11").Value = 5357142.47
$tmeAll.Range("F11").Value = 2511220.78
$tmeAll.Range("G11").Value = 1047985.79
$tmeAll.Range("H11").Value = 2064153.78
$tmeAll.Range("I11").Value = 460510.3
$tmeAll.Range("J11").Value = 576204.22

# ---------------------------------------------------------------------------
# 3. TME_PROV — clear incorrect "SUB IPA 1" tags for Bend Medical Center
# (rows 23-25) and Coastal Health Alliance (rows 37-40)
# ---------------------------------------------------------------------------
$tmeProv = $wb.Worksheets.Item("3. TME_PROV")

$tmeProv.Range("D23").Value = ""
$tmeProv.Range("D24").Value = ""
$tmeProv.Range("D25").Value = ""

$tmeProv.Range("D37").Value = ""
$tmeProv.Range("D38").Value = ""
$tmeProv.Range("D39").Value = ""
$tmeProv.Range("D40").Value = ""

# ---------------------------------------------------------------------------
# 6. RX_MED_PROV — mirror the new TINs from PROV_ID (column C, rows 10-43)
# Force text format first ("@") so the 9-digit TIN strings don't get
# auto-coerced into numbers (these cells carry no explicit style in the
# source file, unlike PROV_ID's column C which is already Text-formatted).
# ---------------------------------------------------------------------------
$rxMedProv = $wb.Worksheets.Item("6. RX_MED_PROV")
$rxMedProv.Range("C10:C43").NumberFormat = "@"

$rxMedProv.Range("C10").Value = "786579303"
$rxMedProv.Range("C11").Value = "786579303"

$rxMedProv.Range("C12").Value = "896233790"
$rxMedProv.Range("C13").Value = "896233790"
$rxMedProv.Range("C14").Value = "896233790"

$rxMedProv.Range("C15").Value = "339670711"
$rxMedProv.Range("C16").Value = "339670711"
$rxMedProv.Range("C17").Value = "339670711"

$rxMedProv.Range("C18").Value = "210053353"
$rxMedProv.Range("C19").Value = "210053353"
$rxMedProv.Range("C20").Value = "210053353"
$rxMedProv.Range("C21").Value = "210053353"

$rxMedProv.Range("C22").Value = "685582861"
$rxMedProv.Range("C23").Value = "685582861"
$rxMedProv.Range("C24").Value = "685582861"
$rxMedProv.Range("C25").Value = "685582861"

$rxMedProv.Range("C26").Value = "553035110"
$rxMedProv.Range("C27").Value = "553035110"
$rxMedProv.Range("C28").Value = "553035110"
$rxMedProv.Range("C29").Value = "553035110"

$rxMedProv.Range("C30").Value = "200604502"
$rxMedProv.Range("C31").Value = "200604502"
$rxMedProv.Range("C32").Value = "200604502"

$rxMedProv.Range("C33").Value = "642621108"
$rxMedProv.Range("C34").Value = "642621108"
$rxMedProv.Range("C35").Value = "642621108"
$rxMedProv.Range("C36").Value = "642621108"

$rxMedProv.Range("C37").Value = "702632297"
$rxMedProv.Range("C38").Value = "702632297"
$rxMedProv.Range("C39").Value = "702632297"
$rxMedProv.Range("C40").Value = "702632297"

$rxMedProv.Range("C41").Value = "797808098"
$rxMedProv.Range("C42").Value = "797808098"
$rxMedProv.Range("C43").Value = "797808098"
